# Add a new "Start Here" task (with a link to the ArcGIS draw-surface
# help topic) at the very top of the task list, followed by a blank
# spacer paragraph, leaving the rest of the document untouched.

$d = $word.ActiveDocument

# Insert two brand-new, run-less paragraphs at the start of the document:
#   1) will become "Start Here: <hyperlink>"
#   2) stays blank, as a spacer before the existing "Measure" task
$insertionPoint = $d.Range(0, 0)
$newParagraphs = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' + `
                  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$insertionPoint.InsertXML($newParagraphs)

# The (still empty) first paragraph becomes the hyperlink run. Hyperlinks.Add
# always places the new hyperlink at the start of its paragraph, so add it
# while the paragraph is empty.
$linkPara = $d.Paragraphs(1).Range
[void]$d.Hyperlinks.Add($linkPara, "http://resources.arcgis.com/en/help/silverlight-api/concepts/index.html#/Using_a_Draw_surface/01660000000s000000/")

# Prefix the paragraph with its lead-in text so the run order ends up as
# "Start Here: " followed by the hyperlink.
$leadIn = $d.Paragraphs(1).Range
$leadIn.Collapse(1)
$leadIn.InsertBefore("Start Here: ")
